# Auto-generated edit script: apply per-cell numeric updates to the
# Masamune Profits workbook (Leve profit-tracking sheets ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 23879.889
$ws.Range("J16").Value = 23879.889
$ws.Range("L16").Value = 23879.889
$ws.Range("N16").Value = -24339.889
$ws.Range("H48").Value = 1136.3636
$ws.Range("J48").Value = 1750
$ws.Range("L48").Value = 5250
$ws.Range("N48").Value = -5834
$ws.Range("H56").Value = 1136.3636
$ws.Range("J56").Value = 1750
$ws.Range("L56").Value = 5250
$ws.Range("N56").Value = -6318
$ws.Range("H63").Value = 37500
$ws.Range("J63").Value = 37500
$ws.Range("L63").Value = 37500
$ws.Range("N63").Value = -38748
$ws.Range("H66").Value = 37500
$ws.Range("J66").Value = 37500
$ws.Range("L66").Value = 112500
$ws.Range("N66").Value = -118740
$ws.Range("H129").Value = 1191.2
$ws.Range("J129").Value = 864.2
$ws.Range("L129").Value = 2592.6
$ws.Range("N129").Value = -12592.6

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2935.7144
$ws.Range("I2").Value = 2258.3333
$ws.Range("K2").Value = 2258.3333
$ws.Range("M2").Value = -2145.3333
$ws.Range("H4").Value = 335
$ws.Range("I4").Value = 270
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 270
$ws.Range("L4").Value = 400
$ws.Range("M4").Value = -154
$ws.Range("N4").Value = -632
$ws.Range("H5").Value = 349.81818
$ws.Range("I5").Value = 231
$ws.Range("J5").Value = 666.6667
$ws.Range("K5").Value = 231
$ws.Range("L5").Value = 666.6667
$ws.Range("M5").Value = -119
$ws.Range("N5").Value = -890.6667
$ws.Range("H35").Value = 3833.3333
$ws.Range("I35").Value = 3833.3333
$ws.Range("K35").Value = 3833.3333
$ws.Range("M35").Value = -3427.3333
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H61").Value = 2264.3667
$ws.Range("I61").Value = 1762.421
$ws.Range("J61").Value = 3131.3635
$ws.Range("K61").Value = 1762.421
$ws.Range("L61").Value = 3131.3635
$ws.Range("M61").Value = -1550.421
$ws.Range("N61").Value = -3555.3635
$ws.Range("H74").Value = 1344.5769
$ws.Range("I74").Value = 1062.8959
$ws.Range("J74").Value = 4724.75
$ws.Range("K74").Value = 1062.8959
$ws.Range("L74").Value = 4724.75
$ws.Range("M74").Value = -188.8959
$ws.Range("N74").Value = -6472.75
$ws.Range("H77").Value = 1344.5769
$ws.Range("I77").Value = 1062.8959
$ws.Range("J77").Value = 4724.75
$ws.Range("K77").Value = 5314.479499999999
$ws.Range("L77").Value = 23623.75
$ws.Range("M77").Value = -946.4794999999995
$ws.Range("N77").Value = -32359.75
$ws.Range("H116").Value = 2935.7144
$ws.Range("I116").Value = 2258.3333
$ws.Range("K116").Value = 2258.3333
$ws.Range("M116").Value = 35.66670000000022
$ws.Range("H122").Value = 1976.8636
$ws.Range("I122").Value = 1935.3529
$ws.Range("J122").Value = 2118
$ws.Range("K122").Value = 5806.0587
$ws.Range("L122").Value = 6354
$ws.Range("M122").Value = -3356.0587
$ws.Range("N122").Value = -11254
$ws.Range("H123").Value = 49213.25
$ws.Range("J123").Value = 49213.25
$ws.Range("L123").Value = 49213.25
$ws.Range("N123").Value = -59013.25
$ws.Range("H132").Value = 19233700
$ws.Range("I132").Value = 31251964
$ws.Range("J132").Value = 4479.4
$ws.Range("K132").Value = 93755892
$ws.Range("L132").Value = 13438.2
$ws.Range("M132").Value = -93753362
$ws.Range("N132").Value = -18498.2
$ws.Range("H136").Value = 2264.3667
$ws.Range("I136").Value = 1762.421
$ws.Range("J136").Value = 3131.3635
$ws.Range("K136").Value = 5287.263
$ws.Range("L136").Value = 9394.0905
$ws.Range("M136").Value = -2737.263
$ws.Range("N136").Value = -14494.0905

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2935.7144
$ws.Range("I3").Value = 2258.3333
$ws.Range("K3").Value = 2258.3333
$ws.Range("M3").Value = -2144.3333
$ws.Range("H4").Value = 349.81818
$ws.Range("I4").Value = 231
$ws.Range("J4").Value = 666.6667
$ws.Range("K4").Value = 231
$ws.Range("L4").Value = 666.6667
$ws.Range("M4").Value = -116
$ws.Range("N4").Value = -896.6667
$ws.Range("H13").Value = 30000
$ws.Range("J13").Value = 30000
$ws.Range("L13").Value = 30000
$ws.Range("N13").Value = -30336
$ws.Range("H20").Value = 4758.636
$ws.Range("I20").Value = 1186.6923
$ws.Range("J20").Value = 7080.4
$ws.Range("K20").Value = 1186.6923
$ws.Range("L20").Value = 7080.4
$ws.Range("M20").Value = -939.6922999999999
$ws.Range("N20").Value = -7574.4
$ws.Range("H36").Value = 3010.25
$ws.Range("I36").Value = 1000
$ws.Range("J36").Value = 3680.3333
$ws.Range("K36").Value = 1000
$ws.Range("L36").Value = 3680.3333
$ws.Range("M36").Value = -466
$ws.Range("N36").Value = -4748.3333
$ws.Range("H107").Value = 2468
$ws.Range("I107").Value = 2468
$ws.Range("K107").Value = 2468
$ws.Range("M107").Value = -548
$ws.Range("H112").Value = 59999
$ws.Range("J112").Value = 59999
$ws.Range("L112").Value = 59999
$ws.Range("N112").Value = -62953
$ws.Range("H134").Value = 2589.7222
$ws.Range("I134").Value = 1943.25
$ws.Range("K134").Value = 5829.75
$ws.Range("M134").Value = -3294.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 19088
$ws.Range("J9").Value = 19088
$ws.Range("L9").Value = 19088
$ws.Range("N9").Value = -19424
$ws.Range("H20").Value = 70000
$ws.Range("J20").Value = 70000
$ws.Range("L20").Value = 70000
$ws.Range("N20").Value = -70472
$ws.Range("H30").Value = 70000
$ws.Range("J30").Value = 70000
$ws.Range("L30").Value = 70000
$ws.Range("N30").Value = -70182
$ws.Range("H47").Value = 26000
$ws.Range("I47").Value = 19000
$ws.Range("J47").Value = 33000
$ws.Range("K47").Value = 19000
$ws.Range("L47").Value = 33000
$ws.Range("M47").Value = -18434
$ws.Range("N47").Value = -34132
$ws.Range("H105").Value = 2418.9473
$ws.Range("I105").Value = 2810.6428
$ws.Range("J105").Value = 1322.2
$ws.Range("K105").Value = 2810.6428
$ws.Range("L105").Value = 1322.2
$ws.Range("M105").Value = -1063.6428
$ws.Range("N105").Value = -4816.2
$ws.Range("H128").Value = 70000
$ws.Range("J128").Value = 70000
$ws.Range("L128").Value = 70000
$ws.Range("N128").Value = -79960

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 8800
$ws.Range("J74").Value = 8800
$ws.Range("L74").Value = 26400
$ws.Range("N74").Value = -28522
$ws.Range("H77").Value = 8800
$ws.Range("J77").Value = 8800
$ws.Range("L77").Value = 79200
$ws.Range("N77").Value = -89808
$ws.Range("H88").Value = 6078.7827
$ws.Range("J88").Value = 6078.7827
$ws.Range("L88").Value = 18236.3481
$ws.Range("N88").Value = -19092.3481
$ws.Range("H91").Value = 6078.7827
$ws.Range("J91").Value = 6078.7827
$ws.Range("L91").Value = 18236.3481
$ws.Range("N91").Value = -21200.3481
$ws.Range("H107").Value = 2668.043
$ws.Range("I107").Value = 3962.111
$ws.Range("J107").Value = 2138.6516
$ws.Range("K107").Value = 11886.333
$ws.Range("L107").Value = 6415.9548
$ws.Range("M107").Value = -9966.332999999999
$ws.Range("N107").Value = -10255.9548
$ws.Range("H131").Value = 922.76
$ws.Range("I131").Value = 677
$ws.Range("K131").Value = 2031
$ws.Range("M131").Value = 3009
$ws.Range("H137").Value = 9478.866
$ws.Range("J137").Value = 12703.667
$ws.Range("L137").Value = 38111.001
$ws.Range("N137").Value = -48311.001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5966.615
$ws.Range("I80").Value = 5153.2354
$ws.Range("J80").Value = 6595.136
$ws.Range("K80").Value = 5153.2354
$ws.Range("L80").Value = 6595.136
$ws.Range("M80").Value = -4155.2354
$ws.Range("N80").Value = -8591.136
$ws.Range("H83").Value = 5966.615
$ws.Range("I83").Value = 5153.2354
$ws.Range("J83").Value = 6595.136
$ws.Range("K83").Value = 25766.177
$ws.Range("L83").Value = 32975.68
$ws.Range("M83").Value = -20774.177
$ws.Range("N83").Value = -42959.68
$ws.Range("H102").Value = 1977.75
$ws.Range("I102").Value = 1244.4
$ws.Range("J102").Value = 3200
$ws.Range("K102").Value = 1244.4
$ws.Range("L102").Value = 3200
$ws.Range("M102").Value = 377.5999999999999
$ws.Range("N102").Value = -6444
$ws.Range("H132").Value = 2181.0334
$ws.Range("I132").Value = 1783.7391
$ws.Range("J132").Value = 3486.4285
$ws.Range("K132").Value = 5351.2173
$ws.Range("L132").Value = 10459.2855
$ws.Range("M132").Value = -2821.2173
$ws.Range("N132").Value = -15519.2855

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 43794.332
$ws.Range("I122").Value = 49736.383
$ws.Range("K122").Value = 149209.149
$ws.Range("M122").Value = -146759.149
$ws.Range("H132").Value = 4260.7144
$ws.Range("I132").Value = 3057.25
$ws.Range("J132").Value = 5865.3335
$ws.Range("K132").Value = 9171.75
$ws.Range("L132").Value = 17596.0005
$ws.Range("M132").Value = -6641.75
$ws.Range("N132").Value = -22656.0005

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 34994.168
$ws.Range("J123").Value = 34994.168
$ws.Range("L123").Value = 34994.168
$ws.Range("N123").Value = -44794.168
$ws.Range("H132").Value = 16130862
$ws.Range("I132").Value = 1586.0625
$ws.Range("J132").Value = 33335424
$ws.Range("K132").Value = 4758.1875
$ws.Range("L132").Value = 100006272
$ws.Range("M132").Value = -2228.1875
$ws.Range("N132").Value = -100011332
